$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '42.098.41'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.491.32'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '312.92'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '95.94'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -3.14%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.558'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -1.78%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.514'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '34.55'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -2.88%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.0790'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '7.09'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -3.45%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '2.879.42'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '2.482.11'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("E16").Value = '  -6.42%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.795'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -4.78%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '42.139.58'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '6.43'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -5.42%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '0.0₃0926'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -2.16%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '11.82'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -2.86%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '69.27'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '238.91'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.83'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -2.57%  '
$ws.Range("E25").Value = '  -4.26%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '24.97'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -5.25%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '2.23'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("E29").Value = '  -2.93%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '37.11'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -6.74%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '155.54'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '5.73'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '2.65'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -6.59%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '2.62'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.0769'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -3.43%  '
$ws.Range("E36").Value = '  -3.57%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '17.28'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -4.29%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '1.91'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -5.93%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.105'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -4.63%  '
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '4.08'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -2.97%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '21.34'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -2.02%  '
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '2.010.39'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +2.79%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.0289'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '3.13'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -4.81%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '8.72'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '2.747.04'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '78.25'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -3.11%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '71.07'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -1.42%  '
$ws.Range("E51").Value = '  -4.07%  '
